$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2441.2
$ws.Range("J17").Value = 1995.8334
$ws.Range("L17").Value = 5987.5002
$ws.Range("N17").Value = -6323.5002

$ws.Range("H33").Value = 407.66666
$ws.Range("I33").Value = 200
$ws.Range("J33").Value = 449.2
$ws.Range("K33").Value = 200
$ws.Range("L33").Value = 449.2
$ws.Range("M33").Value = 29
$ws.Range("N33").Value = -907.2

$ws.Range("H96").Value = 3227.25
$ws.Range("J96").Value = 3254.5
$ws.Range("L96").Value = 9763.5
$ws.Range("N96").Value = -12509.5

$ws.Range("H100").Value = 2193.889
$ws.Range("I100").Value = 1820.8572
$ws.Range("J100").Value = 3499.5
$ws.Range("K100").Value = 1820.8572
$ws.Range("L100").Value = 3499.5
$ws.Range("M100").Value = -1279.8572
$ws.Range("N100").Value = -4581.5

$ws.Range("H107").Value = 378.8
$ws.Range("I107").Value = 173.75
$ws.Range("K107").Value = 173.75
$ws.Range("M107").Value = 1746.25

$ws.Range("H116").Value = 7789.636
$ws.Range("I116").Value = 7633.3335
$ws.Range("J116").Value = 7848.25
$ws.Range("K116").Value = 7633.3335
$ws.Range("L116").Value = 7848.25
$ws.Range("M116").Value = -4191.3335
$ws.Range("N116").Value = -14732.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 395.1
$ws.Range("I88").Value = 395.1
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 395.1
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = 10.89999999999998
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 395.1
$ws.Range("I91").Value = 395.1
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 395.1
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = 1008.9
$ws.Range("N91").ClearContents()

$ws.Range("H122").Value = 557452.25
$ws.Range("I122").Value = 626540.0600000001
$ws.Range("K122").Value = 1879620.18
$ws.Range("M122").Value = -1877170.18

$ws.Range("H132").Value = 2425.9614
$ws.Range("I132").Value = 2425.9614
$ws.Range("K132").Value = 7277.8842
$ws.Range("M132").Value = -4747.8842

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1231.8889
$ws.Range("I20").Value = 1288.4
$ws.Range("K20").Value = 1288.4
$ws.Range("M20").Value = -1041.4

$ws.Range("H80").Value = 209.7
$ws.Range("J80").Value = 40
$ws.Range("L80").Value = 40
$ws.Range("N80").Value = -2036

$ws.Range("H83").Value = 209.7
$ws.Range("J83").Value = 40
$ws.Range("L83").Value = 200
$ws.Range("N83").Value = -10184

$ws.Range("H86").Value = 3994.3635
$ws.Range("I86").Value = 3881.5
$ws.Range("K86").Value = 3881.5
$ws.Range("M86").Value = -2758.5

$ws.Range("H89").Value = 3994.3635
$ws.Range("I89").Value = 3881.5
$ws.Range("K89").Value = 19407.5
$ws.Range("M89").Value = -13791.5

$ws.Range("H134").Value = 2350.0557
$ws.Range("I134").Value = 2086.8
$ws.Range("K134").Value = 6260.400000000001
$ws.Range("M134").Value = -3725.400000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 859.6
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws.Range("H62").Value = 82348.60000000001
$ws.Range("I62").Value = 2996
$ws.Range("J62").Value = 102186.75
$ws.Range("K62").Value = 2996
$ws.Range("L62").Value = 102186.75
$ws.Range("M62").Value = -2372
$ws.Range("N62").Value = -103434.75

$ws.Range("H65").Value = 82348.60000000001
$ws.Range("I65").Value = 2996
$ws.Range("J65").Value = 102186.75
$ws.Range("K65").Value = 14980
$ws.Range("L65").Value = 510933.75
$ws.Range("M65").Value = -11860
$ws.Range("N65").Value = -517173.75

$ws.Range("H86").Value = 11738.2
$ws.Range("J86").Value = 11948.75
$ws.Range("L86").Value = 11948.75
$ws.Range("N86").Value = -14194.75

$ws.Range("H89").Value = 11738.2
$ws.Range("J89").Value = 11948.75
$ws.Range("L89").Value = 59743.75
$ws.Range("N89").Value = -70975.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 147.5
$ws.Range("I8").Value = 147.5
$ws.Range("K8").Value = 442.5
$ws.Range("M8").Value = -303.5

$ws.Range("H14").Value = 767.8461
$ws.Range("I14").Value = 767.8461
$ws.Range("K14").Value = 2303.5383
$ws.Range("M14").Value = -2130.5383

$ws.Range("H62").Value = 1349.5
$ws.Range("I62").Value = 1799
$ws.Range("K62").Value = 5397
$ws.Range("M62").Value = -4711

$ws.Range("H65").Value = 1349.5
$ws.Range("I65").Value = 1799
$ws.Range("K65").Value = 16191
$ws.Range("M65").Value = -12759

$ws.Range("H121").Value = 13828.875
$ws.Range("J121").Value = 15661.857
$ws.Range("L121").Value = 46985.571
$ws.Range("N121").Value = -49605.571

$ws.Range("H131").Value = 1731.625
$ws.Range("I131").Value = 710
$ws.Range("K131").Value = 2130
$ws.Range("M131").Value = 2910

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 20000000
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 20000000
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 20000000
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -20000278

$ws.Range("H52").Value = 19990
$ws.Range("J52").Value = 19990
$ws.Range("L52").Value = 19990
$ws.Range("N52").Value = -20508

$ws.Range("H97").Value = 1203.0625
$ws.Range("I97").Value = 900.8889
$ws.Range("K97").Value = 900.8889
$ws.Range("M97").Value = -404.8889

$ws.Range("H102").Value = 1697.8636
$ws.Range("I102").Value = 416.86667
$ws.Range("J102").Value = 4442.857
$ws.Range("K102").Value = 416.86667
$ws.Range("L102").Value = 416.86667
$ws.Range("M102").Value = 1205.13333
$ws.Range("N102").Value = -7686.857

$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -34900

$ws.Range("H132").Value = 1864.7059
$ws.Range("I132").Value = 1470.4667
$ws.Range("K132").Value = 4411.4001
$ws.Range("M132").Value = -1881.4001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()

$ws.Range("H122").Value = 11129.934
$ws.Range("I122").Value = 10919.385
$ws.Range("K122").Value = 32758.155
$ws.Range("M122").Value = -30308.155

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 1000000
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

$ws.Range("H55").Value = 8000
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H107").Value = 1003.36365
$ws.Range("I107").Value = 428
$ws.Range("K107").Value = 1284
$ws.Range("M107").Value = 636
